# department_200219165142.xlsx -- "update check staff and import excel"
#
# The staff-import sample sheet gets: the sta_sex column switched from a
# free-text "Nam" label to the coded numeric value used by the importer,
# an e-mail hyperlink added to the sample email cell, a second sample row
# (duplicate of the first) appended with its own hyperlink, the
# sta_identity_card_date column widened to fit its header, and the saved
# cursor position / sheet print orientation refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sta_sex (J2): was the shared text string "Nam", now the numeric
#     status code "1" the importer actually expects.
$ws.Range("J2").Value = 1

# --- sta_email (F2): turn the sample address into a real mailto link
#     (Excel auto-creates the built-in "Hyperlink" cell style for this).
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:trung.tran@boot.ai")

# --- Duplicate the sample data row so the sheet demonstrates a second
#     import record, preserving values/number formats from row 2.
$ws.Range("A2:M2").Copy()
$ws.Range("A3:M3").Insert()

# New row's email cell gets its own hyperlink too.
$ws.Hyperlinks.Add($ws.Range("F3"), "mailto:trung.tran@boot.ai")

# Both sample rows get a slightly taller row height.
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 15

# Widen the sta_identity_card_date column (M) so the header fits.
$ws.Columns.Item(13).ColumnWidth = 19

# Refresh the saved selection / cursor position.
$ws.Range("F15").Select() | Out-Null

# Printed page orientation for the sheet.
$ws.PageSetup.Orientation = 1
